$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text-formatted numbers (e.g. "300.80", "0.800") that
# must keep their original string formatting (trailing zeros, thousands dots)
# instead of being auto-coerced into numeric values by Excel. Only the cells
# whose new value would otherwise lose a significant trailing zero need the
# explicit text format applied before the write.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"

$ws.Range("D2").Value = "42.713.88"
$ws.Range("E2").Value = "  +1.23%  "

$ws.Range("D3").Value = "2.294.54"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "300.80"
$ws.Range("E5").Value = "  +0.38%  "

$ws.Range("D6").Value = "98.81"
$ws.Range("E6").Value = "  +3.51%  "

$ws.Range("E7").Value = "  +1.09%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "0.509"
$ws.Range("E9").Value = "  +4.12%  "

$ws.Range("D10").Value = "35.67"
$ws.Range("E10").Value = "  +7.86%  "

$ws.Range("D11").Value = "0.0787"
$ws.Range("E11").Value = "  -0.04%  "

$ws.Range("E12").Value = "  +2.28%  "

$ws.Range("D13").Value = "17.95"
$ws.Range("E13").Value = "  +12.06%  "

$ws.Range("D14").Value = "6.79"
$ws.Range("E14").Value = "  +1.78%  "

$ws.Range("D15").Value = "2.654.73"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").Value = "2.301.15"
$ws.Range("E16").Value = "  +1.78%  "

$ws.Range("D17").Value = "0.800"
$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("D18").Value = "42.629.80"
$ws.Range("E18").Value = "  +1.13%  "

$ws.Range("D19").Value = "12.39"
$ws.Range("E19").Value = "  +6.12%  "

$ws.Range("D20").Value = "6.16"
$ws.Range("E20").Value = "  +2.99%  "

$ws.Range("D21").Value = "0.0₃0896"
$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").Value = "67.65"
$ws.Range("E22").Value = "  +1.98%  "

$ws.Range("D23").Value = "235.18"
$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  +12.80%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Value = "2.44"
$ws.Range("E26").Value = "  -0.40%  "

$ws.Range("D27").Value = "24.45"
$ws.Range("E27").Value = "  +3.09%  "

$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  +6.44%  "

$ws.Range("D29").Value = "167.47"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").Value = "34.09"
$ws.Range("E30").Value = "  +1.78%  "

$ws.Range("D31").Value = "9.11"
$ws.Range("E31").Value = "  +0.08%  "

$ws.Range("E32").Value = "  +0.03%  "

$ws.Range("D33").Value = "4.96"
$ws.Range("E33").Value = "  +1.30%  "

$ws.Range("D34").Value = "4.58"
$ws.Range("E34").Value = "  -2.69%  "

$ws.Range("D35").Value = "17.15"
$ws.Range("E35").Value = "  +2.63%  "

$ws.Range("D36").Value = "2.41"
$ws.Range("E36").Value = "  +3.38%  "

$ws.Range("D37").Value = "0.0685"
$ws.Range("E37").Value = "  -0.51%  "

$ws.Range("D38").Value = "0.101"
$ws.Range("E38").Value = "  +2.84%  "

$ws.Range("D39").Value = "2.82"
$ws.Range("E39").Value = "  +0.69%  "

$ws.Range("E40").Value = "  +2.58%  "

$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("D42").Value = "1.977.89"
$ws.Range("E42").Value = "  +0.90%  "

$ws.Range("E43").Value = "  +2.98%  "

$ws.Range("D44").Value = "2.20"
$ws.Range("E44").Value = "  -5.18%  "

$ws.Range("D45").Value = "10.10"
$ws.Range("E45").Value = "  +5.30%  "

$ws.Range("D46").Value = "17.46"
$ws.Range("E46").Value = "  -1.02%  "

$ws.Range("D47").Value = "2.86"
$ws.Range("E47").Value = "  +3.19%  "

$ws.Range("D48").Value = "55.41"
$ws.Range("E48").Value = "  +6.31%  "

$ws.Range("D49").Value = "2.523.78"
$ws.Range("E49").Value = "  +1.06%  "

$ws.Range("E50").Value = "  +2.94%  "

$ws.Range("D51").Value = "4.49"
$ws.Range("E51").Value = "  +0.41%  "

